$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Negate values in B3:B15 and C3:C15
for ($r = 3; $r -le 15; $r++) {
    $bCell = $ws1.Cells.Item($r, 2)
    $bCell.Value2 = -($bCell.Value2)

    $cCell = $ws1.Cells.Item($r, 3)
    $cCell.Value2 = -($cCell.Value2)
}

$ws1.Range("D1").Select()

$ws1.Activate()

$wb.Save()
